$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 798
$ws.Range("I41").Value = 755.4286
$ws.Range("J41").Value = 872.5
$ws.Range("K41").Value = 755.4286
$ws.Range("L41").Value = 872.5
$ws.Range("M41").Value = -315.4286
$ws.Range("N41").Value = -1752.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6672776.5
$ws.Range("I32").Value = 5090.256
$ws.Range("J32").Value = 47631420
$ws.Range("K32").Value = 5090.256
$ws.Range("L32").Value = 47631420
$ws.Range("M32").Value = -4803.256
$ws.Range("N32").Value = -47631994

$ws.Range("H61").Value = 12503179
$ws.Range("I61").Value = 20835882
$ws.Range("J61").Value = 4123.5
$ws.Range("K61").Value = 20835882
$ws.Range("L61").Value = 4123.5
$ws.Range("M61").Value = -20835670
$ws.Range("N61").Value = -4547.5

$ws.Range("H132").Value = 3679419.5
$ws.Range("I132").Value = 2860.4614
$ws.Range("J132").Value = 19611176
$ws.Range("K132").Value = 8581.3842
$ws.Range("L132").Value = 58833528
$ws.Range("M132").Value = -6051.3842
$ws.Range("N132").Value = -58838588

$ws.Range("H136").Value = 12503179
$ws.Range("I136").Value = 20835882
$ws.Range("J136").Value = 4123.5
$ws.Range("K136").Value = 62507646
$ws.Range("L136").Value = 12370.5
$ws.Range("M136").Value = -62505096
$ws.Range("N136").Value = -17470.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 15437.714
$ws.Range("I82").Value = 7289.25
$ws.Range("J82").Value = 26302.334
$ws.Range("K82").Value = 7289.25
$ws.Range("L82").Value = 26302.334
$ws.Range("M82").Value = -6906.25
$ws.Range("N82").Value = -27068.334

$ws.Range("H85").Value = 15437.714
$ws.Range("I85").Value = 7289.25
$ws.Range("J85").Value = 26302.334
$ws.Range("K85").Value = 7289.25
$ws.Range("L85").Value = 26302.334
$ws.Range("M85").Value = -5963.25
$ws.Range("N85").Value = -28954.334

$ws.Range("H134").Value = 13905364
$ws.Range("I134").Value = 5200
$ws.Range("J134").Value = 27805528
$ws.Range("K134").Value = 15600
$ws.Range("L134").Value = 83416584
$ws.Range("M134").Value = -13065
$ws.Range("N134").Value = -83421654

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 300
$ws.Range("I2").Value = 300
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 300
$ws.Range("M2").Value = -187
$ws.Range("N2").Value = -526

$ws.Range("H31").Value = 995.3125
$ws.Range("I31").Value = 787.931
$ws.Range("K31").Value = 787.931
$ws.Range("M31").Value = -492.931

$ws.Range("H34").Value = 995.3125
$ws.Range("I34").Value = 787.931
$ws.Range("K34").Value = 787.931
$ws.Range("M34").Value = -585.931

$ws.Range("H86").Value = 21785470
$ws.Range("I86").Value = 62564964
$ws.Range("J86").Value = 36407.2
$ws.Range("K86").Value = 62564964
$ws.Range("L86").Value = 36407.2
$ws.Range("M86").Value = -62563841
$ws.Range("N86").Value = -38653.2

$ws.Range("H89").Value = 21785470
$ws.Range("I89").Value = 62564964
$ws.Range("J89").Value = 36407.2
$ws.Range("K89").Value = 312824820
$ws.Range("L89").Value = 182036
$ws.Range("M89").Value = -312819204
$ws.Range("N89").Value = -193268

$ws.Range("H132").Value = 13335508
$ws.Range("I132").Value = 1759.875
$ws.Range("K132").Value = 5279.625
$ws.Range("M132").Value = -2749.625

$ws.Range("H134").Value = 1510.6
$ws.Range("I134").Value = 1172.8125
$ws.Range("J134").Value = 2111.111
$ws.Range("K134").Value = 3518.4375
$ws.Range("L134").Value = 6333.333
$ws.Range("M134").Value = -983.4375
$ws.Range("N134").Value = -11403.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 142914.14
$ws.Range("I4").Value = 142914.14
$ws.Range("K4").Value = 428742.42
$ws.Range("M4").Value = -428630.42

$ws.Range("H92").Value = 5343.1816
$ws.Range("I92").Value = 636.8570999999999
$ws.Range("J92").Value = 7539.467
$ws.Range("K92").Value = 1910.5713
$ws.Range("L92").Value = 22618.401
$ws.Range("M92").Value = -662.5712999999998
$ws.Range("N92").Value = -25114.401

$ws.Range("H113").Value = 4472157.5
$ws.Range("I113").Value = 3788637.2
$ws.Range("J113").Value = 5263601.5
$ws.Range("K113").Value = 11365911.6
$ws.Range("L113").Value = 15790804.5
$ws.Range("M113").Value = -11363741.6
$ws.Range("N113").Value = -15795144.5

$ws.Range("H131").Value = 887.65
$ws.Range("I131").Value = 712.5
$ws.Range("J131").Value = 894.94794
$ws.Range("K131").Value = 2137.5
$ws.Range("L131").Value = 2684.84382
$ws.Range("M131").Value = 2902.5
$ws.Range("N131").Value = -12764.84382

$ws.Range("H132").Value = 8929.643
$ws.Range("J132").Value = 9501.154
$ws.Range("L132").Value = 85510.386
$ws.Range("N132").Value = -90570.386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 3002282.2
$ws.Range("I12").Value = 3002282.2
$ws.Range("K12").Value = 3002282.2
$ws.Range("M12").Value = -3002142.2

$ws.Range("H136").Value = 22056.5
$ws.Range("J136").Value = 22056.5
$ws.Range("L136").Value = 66169.5
$ws.Range("N136").Value = -71269.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 500
$ws.Range("J25").Value = 500
$ws.Range("L25").Value = 500
$ws.Range("N25").Value = -960

$ws.Range("H132").Value = 66685900
$ws.Range("I132").Value = 200003860
$ws.Range("J132").Value = 26921
$ws.Range("K132").Value = 600011580
$ws.Range("L132").Value = 80763
$ws.Range("M132").Value = -600009050
$ws.Range("N132").Value = -85823

$ws.Range("H136").Value = 40819196
$ws.Range("I136").Value = 10207006
$ws.Range("J136").Value = 71431384
$ws.Range("K136").Value = 30621018
$ws.Range("L136").Value = 214294152
$ws.Range("M136").Value = -30618468
$ws.Range("N136").Value = -214299252

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 5000
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5298

$ws.Range("H126").Value = 2987.375
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 3149.8333
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 9449.499899999999
$ws.Range("M126").Value = -5030
$ws.Range("N126").Value = -14389.4999
